{"js": "// Replace \"clarification\" with \"resolution\", and expand\n// \"circumstance among others\" into the longer clause describing how these\n// factors could explain procrastination delays when discounting doesn't.\n\nconst oldTail = \"circumstance among others\";\nconst newTail = \"circumstances which could all explain delays in cases where discounting fails to or is absent altogether\";\n\nconst r1 = context.document.body.search(\"clarification\", { matchCase: true, matchWholeWord: true });\nr1.load(\"text\");\nawait context.sync();\n\nif (r1.items.length > 0) {\n  r1.items[0].insertText(\"resolution\", Word.InsertLocation.replace);\n}\n\nconst r2 = context.document.body.search(oldTail, { matchCase: true });\nr2.load(\"text\");\nawait context.sync();\n\nif (r2.items.length > 0) {\n  r2.items[0].insertText(newTail, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace \"clarification\" with \"resolution\", and expand\n# \"circumstance among others\" into the longer clause describing how these\n# factors could explain procrastination delays when discounting doesn't.\n\n$d = $word.ActiveDocument\n\n$range1 = $d.Content\n$range1.Find.Execute(\n    \"clarification\",  # FindText\n    $false,            # MatchCase\n    $true,             # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    \"resolution\",      # ReplaceWith\n    2                  # Replace (wdReplaceOne)\n) | Out-Null\n\n$range2 = $d.Content\n$range2.Find.Execute(\n    \"circumstance among others\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"circumstances which could all explain delays in cases where discounting fails to or is absent altogether\",\n    2\n) | Out-Null\n"}
